# Rename the "2108 LEAVE CREDITS" sheet to "2018 LEAVE CREDITS"
# (the original sheet name had a typo: 2108 instead of 2018)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2108 LEAVE CREDITS")
$ws.Name = "2018 LEAVE CREDITS"

# The sheet's print-title rows defined name embeds the old sheet name in its
# formula text ('2108 LEAVE CREDITS'!$1:$9) and does not auto-update when the
# sheet is renamed, so refresh it explicitly to reference the new name.
$ws.PageSetup.PrintTitleRows = "'2018 LEAVE CREDITS'!`$1:`$9"
